# Continue the "Excel to mySQL" staging table on Sheet1 with one more
# person's record (row 12): id, first name, city, number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "joma"
$ws.Range("C12").Value = "suli"
$ws.Range("D12").Value = 8948583

# Leave the selection where the user's cursor ended up after entering the
# row (one cell to the right of the last value).
$ws.Range("E12").Select()
